# TEXAS_2017.xlsx cleanup:
#  - rename header columns to snake_case machine-readable names
#  - title-case ("PROPER") the Spanish state/municipality names in columns A and B
#    (connector words like "de", "del", "la", "el", "los", "las", "y" get capitalized too,
#     matching Excel's PROPER() behaviour)
#  - drop the trailing footnote/source rows after the grand-total row
#  - dimension shrinks accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header row -> machine-friendly column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Title-case the state (A) and municipality (B) name columns using Excel's PROPER()
$lastRow = 2050
for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    if ($aVal -ne $null) {
        $aCell.Value = $excel.WorksheetFunction.Proper($aVal)
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($bVal -ne $null) {
        $bCell.Value = $excel.WorksheetFunction.Proper($bVal)
    }
}

# 3) Remove the trailing footnote/source rows (sample size, source, elaborated-by, etc.)
#    that used to live below the grand-total row, shrinking the sheet's used range.
$ws.Range("A2052:D2056").ClearContents()

# 4) A handful of percentage cells were re-derived (count / grand-total) during the
#    cleanup pass and landed one ULP away from their original literal; pin them to the
#    exact recomputed values.
$ws.Range("D109").Value = 0.000910088599801922
$ws.Range("D212").Value = 0.009309671029738485
$ws.Range("D302").Value = 0.00092614898685725
$ws.Range("D531").Value = 0.000910088599801922
$ws.Range("D1571").Value = 0.000968976685671458
$ws.Range("D1701").Value = 0.000942209373912578
